# Add a new fee-earner mapping row for Sandeep Chahil
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new row's values first
$ws.Range("A18").Value = "Sandeep Chahil"
$ws.Range("B18").Value = "Nicola.Daniel@taylorslegal.com"

# Copy formatting from the row above (row 16: A16/B16) down to the new row 18
# (xlPasteFormats = -4122), so the cell styles match the rest of the table
$ws.Range("A16:B16").Copy()
$ws.Range("A18:B18").PasteSpecial(-4122)

# Mirror the author's final selection state
$ws.Range("F18:G19").Select()
